$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RESULTS FRAMEWORK")
$ws.Columns("B").Delete()
